$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SCD0338")

# Update formulas to shift the "month offset" window back by 30 days each
$ws.Range("Q2").Formula = '=TEXT(TODAY()-30,"mmmm")'
$ws.Range("Q3").Formula = '=TEXT(TODAY()-60,"mmmm")'
$ws.Range("Q4").Formula = '=TEXT(TODAY()-30,"mmmm")'

# Update the active selection on the sheet
$ws.Range("Q5").Select()
